$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, border, centered) from P1 into Q1:R1
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)

# Header labels
$ws.Range("Q1").Value = "PD"
$ws.Range("R1").Value = "N2"

# Data for columns Q (PD) and R (N2) per row
$qrData = @{
    2 = @("AGU-E", "Fuera de Poligono OVL")
    3 = @("PUE-B", "Fuera de Poligono OVL")
    4 = @("BLO-?", "Fuera de Poligono OVL")
    5 = @("ATH-?", "Fuera de Poligono OVL")
    6 = @("VCR-I", "Fuera de Poligono OVL")
    7 = @("VCR-F", "Fuera de Poligono OVL")
    8 = @("AGU-O", "Fuera de Poligono OVL")
    9 = @("BLO-A", "Fuera de Poligono OVL")
    10 = @("COG-H", "Fuera de Poligono OVL")
    11 = @("DEV-F", "ARATO-25058.PO.1DEV")
    12 = @("PUE-G", "Fuera de Poligono OVL")
    13 = @("BLO-J", "Fuera de Poligono OVL")
    14 = @("BLO-J", "Fuera de Poligono OVL")
    15 = @("AGU-N", "Fuera de Poligono OVL")
    16 = @("COG-F", "Fuera de Poligono OVL")
    17 = @("DEV-L", "ARATO-25058.PO.2DEV")
    18 = @("VCR-N", "Fuera de Poligono OVL")
    19 = @("DEV-F", "ARATO-25058.PO.1DEV")
    20 = @("VCR-O", "Fuera de Poligono OVL")
    21 = @("PUE-M", "ARATO-25058.PO.1PUE")
    22 = @("PUE-H", "Fuera de Poligono OVL")
    23 = @("ATH-S", "Fuera de Poligono OVL")
    24 = @("PUE-M", "ARATO-25058.PO.1PUE")
    25 = @("PUE-J", "Fuera de Poligono OVL")
    26 = @("PUE-M", "ARATO-25058.PO.1PUE")
    27 = @("COG-A", "Fuera de Poligono OVL")
    28 = @("PUE-F", "Fuera de Poligono OVL")
    29 = @("PUE-I", "Fuera de Poligono OVL")
    30 = @("BLO-?", "Fuera de Poligono OVL")
    31 = @("ALM-N", "Fuera de Poligono OVL")
    32 = @("ALM-O", "Fuera de Poligono OVL")
    33 = @("ATH-J", "Fuera de Poligono OVL")
    34 = @("CON-B", "Fuera de Poligono OVL")
    35 = @("RET-A", "Fuera de Poligono OVL")
    36 = @("BLO-N", "Fuera de Poligono OVL")
    37 = @("COG-?", "Fuera de Poligono OVL")
    38 = @("BLO-F", "Fuera de Poligono OVL")
    39 = @("CLI-O", "Fuera de Poligono OVL")
    40 = @("COG-K", "Fuera de Poligono OVL")
    41 = @("AGU-C", "Fuera de Poligono OVL")
    42 = @("PAV-V", "Fuera de Poligono OVL")
    43 = @("ALM-A", "Fuera de Poligono OVL")
    44 = @("CON-G", "Fuera de Poligono OVL")
    45 = @("PCH-S", "ARATO-25058.PO.2PCH")
    46 = @("CLI-I", "Fuera de Poligono OVL")
    47 = @("CLI-M", "Fuera de Poligono OVL")
    48 = @("CLI-N", "Fuera de Poligono OVL")
    49 = @("DEV-D", "Fuera de Poligono OVL")
    50 = @("AGU-N", "Fuera de Poligono OVL")
    51 = @("BLO-I", "Fuera de Poligono OVL")
    52 = @("VCR-B", "Fuera de Poligono OVL")
    53 = @("PAV-O", "Fuera de Poligono OVL")
    54 = @("ATH-C", "Fuera de Poligono OVL")
    55 = @("VCR-N", "Fuera de Poligono OVL")
    56 = @("NRA-R", "ARATO-25058.PO.2NRA")
    57 = @("CEN-M", "Fuera de Poligono OVL")
    58 = @("AGU-P", "Fuera de Poligono OVL")
    59 = @("NRA-R", "ARATO-25058.PO.2NRA")
    60 = @("CLI-M", "Fuera de Poligono OVL")
    61 = @("PAV-E", "Fuera de Poligono OVL")
    62 = @("VCR-M", "Fuera de Poligono OVL")
    63 = @("PAV-M", "Fuera de Poligono OVL")
    64 = @("CON-I", "Fuera de Poligono OVL")
    65 = @("BLO-F", "Fuera de Poligono OVL")
    66 = @("CON-G", "Fuera de Poligono OVL")
    67 = @("CLI-H", "Fuera de Poligono OVL")
    68 = @("PUE-J", "Fuera de Poligono OVL")
    69 = @("COG-F", "Fuera de Poligono OVL")
    70 = @("PUE-J", "Fuera de Poligono OVL")
    71 = @("PUE-O", "ARATO-25058.PO.1PUE")
    72 = @("COG-O", "Fuera de Poligono OVL")
    73 = @("PPT-H", "Fuera de Poligono OVL")
    74 = @("CLI-N", "Fuera de Poligono OVL")
    75 = @("PPT-K", "Fuera de Poligono OVL")
    76 = @("VCR-O", "Fuera de Poligono OVL")
    77 = @("COG-C", "Fuera de Poligono OVL")
    78 = @("NRA-H", "Fuera de Poligono OVL")
}

foreach ($rowNum in $qrData.Keys) {
    $vals = $qrData[$rowNum]
    $ws.Cells.Item($rowNum, 17).Value = $vals[0]
    $ws.Cells.Item($rowNum, 18).Value = $vals[1]
}

Write-Host "Done updating columns Q and R"